$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.376, 0.155, 0.444, 92, 48, -175),
    @(-0.165, 0.181, 0.662, -61, -13, 21),
    @(0.591, -0.236, 0.193, -61, -25, 156),
    @(0.188, -0.002, 0.628, 30, -59, 158),
    @(-0.463, 0.248, 0.371, -39, -10, -107),
    @(-0.14, 0.057, 0.432, 156, -1, -10),
    @(-0.168, -0.447, 0.469, 51, 0, -134),
    @(-0.333, -0.052, -0.151, 98, -1, -19),
    @(-0.059, -0.026, 0.602, -118, 58, -157),
    @(-0.369, -0.078, 0.312, 31, 18, 173)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$wb.Save()
